$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.850.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.74%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.492.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.56%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.87%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.31%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  +1.35%  '

# Row 9
$ws.Range('E9').Value = '  +1.74%  '

# Row 10
$ws.Range('E10').Value = '  -1.66%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.39'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.04%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.347'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.27%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.931.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.47%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '58.729.71'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.69%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.34%  '

# Row 16
$ws.Range('E16').Value = '  +0.11%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.501.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.52%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.70%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.62%  '

# Row 22
$ws.Range('E22').Value = '  +4.38%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.28%  '

# Row 24
$ws.Range('E24').Value = '  +2.75%  '

# Row 25
$ws.Range('E25').Value = '  +0.39%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.41%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.76%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0759'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.82%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.15%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.87%  '

# Row 31
$ws.Range('E31').Value = '  -0.68%  '

# Row 33
$ws.Range('E33').Value = '  -0.03%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.54%  '

# Row 35
$ws.Range('E35').Value = '  -0.11%  '

# Row 36
$ws.Range('E36').Value = '  +0.04%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.83%  '

# Row 38
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.797'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.22%  '

# Row 39
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.38%  '

# Row 40
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '281.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.48%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.30%  '

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.03%  '

# Row 43
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.603'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.97%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '129.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.09%  '

# Row 45
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.61%  '

# Row 46
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0921'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.16%  '

# Row 47
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0500'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.07%  '

# Row 48
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0217'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.41%  '

# Row 49
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.34%  '

# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.752.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.51%  '

# Row 51
$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.980'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.37%  '
